# Updates the cryptocurrency symbol list (commit: "Updated symbol list on
# Fri Jan  6 13:30:28 UTC 2023 with GitHub Actions").
#
# For each changed row we update the Coin (B), Link (C), Price (D) and
# Volume/1h (E) columns to their new values. D and E hold numeric-looking
# text (e.g. "255.15" or "-0.58%") that must stay stored as plain text
# (matching the source workbook's inlineStr cells), so each of those
# assignments is apostrophe-prefixed to stop Excel's automatic
# number/percentage conversion; the leading apostrophe itself is not
# stored as part of the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'255.15"
$ws.Range("E2").Value = "'-0.58%"
# Row 3
$ws.Range("D3").Value = "'26.39"
$ws.Range("E3").Value = "'-2.60%"
# Row 4
$ws.Range("D4").Value = "'4.642"
$ws.Range("E4").Value = "'-1.43%"
# Row 5
$ws.Range("D5").Value = "'0.05923"
$ws.Range("E5").Value = "'0.10%"
# Row 6
$ws.Range("D6").Value = "'6.630"
$ws.Range("E6").Value = "'-0.20%"
# Row 7
$ws.Range("D7").Value = "'0.8525"
$ws.Range("E7").Value = "'-1.75%"
# Row 8
$ws.Range("D8").Value = "'0.9053"
$ws.Range("E8").Value = "'-4.69%"
# Row 9
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.0006029"
$ws.Range("E9").Value = "'-0.01%"
# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1376"
$ws.Range("E10").Value = "'-2.28%"
# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.04165"
$ws.Range("E11").Value = "'8.73%"
# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.06985"
$ws.Range("E12").Value = "'-1.57%"
# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03032"
$ws.Range("E13").Value = "'-5.76%"
# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09091"
$ws.Range("E14").Value = "'-1.86%"
# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001534"
$ws.Range("E15").Value = "'0.04%"
# Row 16
$ws.Range("D16").Value = "'0.006116"
$ws.Range("E16").Value = "'1.56%"
# Row 17
$ws.Range("D17").Value = "'3.470"
$ws.Range("E17").Value = "'-1.26%"
# Row 18
$ws.Range("D18").Value = "'3.143"
$ws.Range("E18").Value = "'-1.60%"
# Row 19
$ws.Range("E19").Value = "'-1.90%"
# Row 20
$ws.Range("D20").Value = "'0.3022"
$ws.Range("E20").Value = "'-3.59%"
# Row 21
$ws.Range("D21").Value = "'0.1284"
$ws.Range("E21").Value = "'0.05%"
# Row 22
$ws.Range("D22").Value = "'3.857"
$ws.Range("E22").Value = "'-0.60%"
# Row 23
$ws.Range("D23").Value = "'0.04207"
$ws.Range("E23").Value = "'-0.61%"
# Row 24
$ws.Range("E24").Value = "'-0.58%"
# Row 25
$ws.Range("D25").Value = "'0.004681"
$ws.Range("E25").Value = "'8.93%"
# Row 26
$ws.Range("E26").Value = "'-0.02%"
# Row 27
$ws.Range("D27").Value = "'0.0001523"
$ws.Range("E27").Value = "'1.45%"
# Row 40
$ws.Range("D40").Value = "'0.03769"
$ws.Range("E40").Value = "'-1.35%"
# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006248"
$ws.Range("E41").Value = "'58.03%"
# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1094"
$ws.Range("E42").Value = "'-0.75%"
# Row 43
$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'-9.12%"
# Row 44
$ws.Range("D44").Value = "'0.01451"
$ws.Range("E44").Value = "'26.11%"
# Row 45
$ws.Range("E45").Value = "'-6.30%"
# Row 46
$ws.Range("E46").Value = "'-0.01%"
# Row 47
$ws.Range("D47").Value = "'0.03999"
$ws.Range("E47").Value = "'-50.32%"
# Row 48
$ws.Range("E48").Value = "'10,460.84%"
# Row 49
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.01%"
# Row 50
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.01%"
